# Update source/input payroll figures on Sheet1 (row 24: angajat 4; row 25/26:
# company contribution rows) so that all dependent totals (rows 27-29) and the
# summary table (row 39) recalculate to the corrected amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - employee values
$ws.Range("M24").Value = 349.0
$ws.Range("Q24").Value = 1327.0
$ws.Range("S24").Value = 348.0

# Row 25 - CM Societate values
$ws.Range("O25").Value = 5307.0
$ws.Range("R25").Value = 3136.0
$ws.Range("S25").Value = 3136.0

# Row 26 - totals row for this block
$ws.Range("R26").Value = 3484.0
$ws.Range("T26").Value = 3136.0

# Row 39 - standalone summary cells (not formulas, mirror Q27 / S27)
$ws.Range("E39").Value = 2594.0
$ws.Range("L39").Value = 600.0

# Recalculate all formulas (rows 27-29 depend on the values changed above)
$excel.CalculateFullRebuild()

# Scroll back to the top-left of the sheet after the recalculation/reset
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
